$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($rng, $val) {
    # Force the cell to Text format before/after assignment so
    # numeric-looking strings (e.g. "57.26") are NOT auto-converted
    # to numbers, and the cell keeps its original (default) style.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "37.843.32"
$ws.Range("E2").Value = "  +1.59%  "

Set-TextCell $ws.Range("D3") "2.085.40"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("E4").Value = "  +0.05%  "

Set-TextCell $ws.Range("D5") "232.44"
$ws.Range("E5").Value = "  -0.60%  "

Set-TextCell $ws.Range("D6") "0.625"
$ws.Range("E6").Value = "  +0.38%  "

Set-TextCell $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.00%  "

Set-TextCell $ws.Range("D8") "57.26"
$ws.Range("E8").Value = "  +0.83%  "

Set-TextCell $ws.Range("D9") "0.388"
$ws.Range("E9").Value = "  +1.31%  "

Set-TextCell $ws.Range("D10") "0.0780"
$ws.Range("E10").Value = "  +2.05%  "

$ws.Range("E11").Value = "  +2.62%  "

Set-TextCell $ws.Range("D12") "2.382.56"
$ws.Range("E12").Value = "  +0.36%  "

Set-TextCell $ws.Range("D13") "14.43"
$ws.Range("E13").Value = "  -1.53%  "

Set-TextCell $ws.Range("D14") "21.06"
$ws.Range("E14").Value = "  +1.21%  "

Set-TextCell $ws.Range("D15") "0.763"
$ws.Range("E15").Value = "  -1.80%  "

Set-TextCell $ws.Range("D16") "5.26"
$ws.Range("E16").Value = "  +2.71%  "

Set-TextCell $ws.Range("D17") "2.087.91"
$ws.Range("E17").Value = "  +0.87%  "

Set-TextCell $ws.Range("D18") "37.776.42"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("E19").Value = "  -4.29%  "

Set-TextCell $ws.Range("D20") "70.67"
$ws.Range("E20").Value = "  +1.73%  "

Set-TextCell $ws.Range("D21") "0.0₃0821"
$ws.Range("E21").Value = "  +1.01%  "

Set-TextCell $ws.Range("D22") "228.14"
$ws.Range("E22").Value = "  +0.70%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("E24").Value = "  -1.03%  "

$ws.Range("E25").Value = "  -0.24%  "

Set-TextCell $ws.Range("D26") "170.51"
$ws.Range("E26").Value = "  +2.24%  "

Set-TextCell $ws.Range("D27") "0.140"
$ws.Range("E27").Value = "  +11.17%  "

Set-TextCell $ws.Range("D28") "8.92"
$ws.Range("E28").Value = "  +1.73%  "

Set-TextCell $ws.Range("D29") "1.43"
$ws.Range("E29").Value = "  +0.00%  "

Set-TextCell $ws.Range("D30") "19.44"
$ws.Range("E30").Value = "  +1.94%  "

$ws.Range("E31").Value = "  +0.96%  "

Set-TextCell $ws.Range("D32") "4.62"
$ws.Range("E32").Value = "  +3.53%  "

Set-TextCell $ws.Range("D33") "0.0624"
$ws.Range("E33").Value = "  +1.30%  "

Set-TextCell $ws.Range("D34") "4.59"
$ws.Range("E34").Value = "  +0.30%  "

Set-TextCell $ws.Range("D35") "2.50"
$ws.Range("E35").Value = "  +0.20%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws.Range("D36") "1.83"
$ws.Range("E36").Value = "  +3.41%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Range("D37") "3.40"
$ws.Range("E37").Value = "  +5.26%  "

$ws.Range("E38").Value = "  +0.01%  "

Set-TextCell $ws.Range("D39") "5.42"
$ws.Range("E39").Value = "  -4.53%  "

Set-TextCell $ws.Range("D40") "0.0998"
$ws.Range("E40").Value = "  +6.89%  "

$ws.Range("E41").Value = "  -0.67%  "

Set-TextCell $ws.Range("D42") "97.33"
$ws.Range("E42").Value = "  +1.31%  "

Set-TextCell $ws.Range("D43") "0.0213"
$ws.Range("E43").Value = "  +0.76%  "

Set-TextCell $ws.Range("D44") "1.448.32"
$ws.Range("E44").Value = "  -2.15%  "

Set-TextCell $ws.Range("D45") "1.16"
$ws.Range("E45").Value = "  -0.19%  "

Set-TextCell $ws.Range("D46") "1.05"
$ws.Range("E46").Value = "  +2.92%  "

Set-TextCell $ws.Range("D47") "15.66"
$ws.Range("E47").Value = "  +3.60%  "

Set-TextCell $ws.Range("D48") "4.06"
$ws.Range("E48").Value = "  -8.31%  "

Set-TextCell $ws.Range("D49") "7.39"
$ws.Range("E49").Value = "  +3.33%  "

Set-TextCell $ws.Range("D50") "3.02"
$ws.Range("E50").Value = "  +1.91%  "

Set-TextCell $ws.Range("D51") "2.278.84"
$ws.Range("E51").Value = "  +0.68%  "

